# Update countries & provincias Spain
# Refresh of the "Pais" COVID dashboard: new case totals pushed a handful of
# countries past their neighbours in the (externally sorted) list, and the
# "last updated" timestamp moved from 12:05 to 13:22. Because the sheet is
# kept in descending "Casos totales" order, re-ranking shows up as the
# country label in a row changing while its stats catch up to the row
# above/below; every other row's data is simply refreshed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: country label updates (rows whose rank/identity rotated) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 13:22"
$ws.Range("A49").Value = "Suiza"
$ws.Range("A50").Value = "Guatemala"
$ws.Range("A51").Value = "Costa Rica"
$ws.Range("A133").Value = "Malta"
$ws.Range("A134").Value = "Polinesia Francesa"
$ws.Range("A135").Value = "Congo"
$ws.Range("A136").Value = "Surinam"
$ws.Range("A188").Value = "Liechtenstein"
$ws.Range("A189").Value = "Camboya"
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("A217").Value = "Montserrat"

# --- Numeric data updates (Casos totales/Nuevos/Activos/Recuperados/Criticos/MuertesHoy/Muertes) ---
$ws.Range("B16").Value = 556891
$ws.Range("C16").Value = 6134
$ws.Range("D16").Value = 446685
$ws.Range("E16").Value = 78221
$ws.Range("G16").Value = 335
$ws.Range("H16").Value = 31985

$ws.Range("B20").Value = 404856
$ws.Range("C20").Value = 982
$ws.Range("D20").Value = 310200
$ws.Range("E20").Value = 84611
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 10045

$ws.Range("B21").Value = 396413
$ws.Range("C21").Value = 1586
$ws.Range("D21").Value = 312065
$ws.Range("E21").Value = 78587
$ws.Range("G21").Value = 14
$ws.Range("H21").Value = 5761

$ws.Range("B28").Value = 308682
$ws.Range("C28").Value = 435
$ws.Range("D28").Value = 289171
$ws.Range("E28").Value = 17183
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 2328

$ws.Range("B34").Value = 201032
$ws.Range("C34").Value = 5028
$ws.Range("D34").Value = 144429
$ws.Range("E34").Value = 50358
$ws.Range("G34").Value = 82
$ws.Range("H34").Value = 6245

$ws.Range("B37").Value = 153008
$ws.Range("C37").Value = 4499
$ws.Range("D37").Value = 105488
$ws.Range("E37").Value = 46691
$ws.Range("G37").Value = 17
$ws.Range("H37").Value = 829

$ws.Range("B39").Value = 130711
$ws.Range("C39").Value = 249
$ws.Range("D39").Value = 127599
$ws.Range("E39").Value = 2883
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 229

$ws.Range("B42").Value = 122273
$ws.Range("C42").Value = 1563
$ws.Range("D42").Value = 115068
$ws.Range("E42").Value = 6730
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 475

$ws.Range("B49").Value = 103653
$ws.Range("C49").Value = 6634
$ws.Range("D49").Value = 55700
$ws.Range("E49").Value = 45897
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 2056

$ws.Range("B50").Value = 103172
$ws.Range("D50").Value = 92665
$ws.Range("E50").Value = 6927
$ws.Range("H50").Value = 3580

$ws.Range("B51").Value = 100616
$ws.Range("D51").Value = 61162
$ws.Range("E51").Value = 38203
$ws.Range("G51").Value = 1251
$ws.Range("H51").Value = 2056

$ws.Range("B71").Value = 53384
$ws.Range("C71").Value = 764
$ws.Range("D71").Value = 29619
$ws.Range("E71").Value = 22991
$ws.Range("G71").Value = 6
$ws.Range("H71").Value = 774

$ws.Range("B88").Value = 27484
$ws.Range("C88").Value = 18
$ws.Range("D88").Value = 25169
$ws.Range("E88").Value = 1410

$ws.Range("B92").Value = 24514
$ws.Range("C92").Value = 710
$ws.Range("D92").Value = 15884
$ws.Range("E92").Value = 8416
$ws.Range("G92").Value = 10
$ws.Range("H92").Value = 214

$ws.Range("B101").Value = 15525
$ws.Range("C101").Value = 17
$ws.Range("D101").Value = 14082
$ws.Range("E101").Value = 1122

$ws.Range("E102").Value = 4321
$ws.Range("H102").Value = 353

$ws.Range("B133").Value = 5258
$ws.Range("C133").Value = 121
$ws.Range("D133").Value = 3439
$ws.Range("E133").Value = 1770
$ws.Range("H133").Value = 49

$ws.Range("B134").Value = 5161
$ws.Range("D134").Value = 3536
$ws.Range("E134").Value = 1606
$ws.Range("H134").Value = 19

$ws.Range("B135").Value = 5156
$ws.Range("D135").Value = 3887
$ws.Range("E135").Value = 1177
$ws.Range("H135").Value = 92

$ws.Range("B136").Value = 5154
$ws.Range("D136").Value = 4995
$ws.Range("E136").Value = 50
$ws.Range("H136").Value = 109

$ws.Range("B175").Value = 641
$ws.Range("C175").Value = 11
$ws.Range("D175").Value = 500
$ws.Range("E175").Value = 141

$ws.Range("B188").Value = 324
$ws.Range("C188").Value = 42
$ws.Range("D188").Value = 170
$ws.Range("E188").Value = 153
$ws.Range("H188").Value = 1

$ws.Range("B189").Value = 286
$ws.Range("D189").Value = 280
$ws.Range("E189").Value = 6
$ws.Range("H189").Value = 0

$ws.Range("B194").Value = 153
$ws.Range("C194").Value = 2
$ws.Range("D194").Value = 149
$ws.Range("E194").Value = 4

$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

$ws.Range("D217").Value = 12
$ws.Range("H217").Value = 1
